$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two brand-new rows before the current row 282; this pushes the
# existing rows 282-366 down to 284-368 (and the sheet dimension grows
# from R366 to R368 automatically).
$ws.Rows("282:283").Insert()

# Populate the first new row (becomes row 282) with fresh market data.
$ws.Range("A282").Value = 10
$ws.Range("B282").Value = "Vega Modelo de Temuco"
$ws.Range("C282").Value = "La Araucanía"
$ws.Range("D282").Value = 44463
$ws.Range("E282").Value = 9
$ws.Range("F282").Value = 100112003
$ws.Range("G282").Value = "Ajo"
$ws.Range("H282").Value = "Chino"
$ws.Range("I282").Value = "Primera"
$ws.Range("J282").Value = 140
$ws.Range("K282").Value = 18000
$ws.Range("L282").Value = 18000
$ws.Range("M282").Value = 18000
$ws.Range("N282").Value = "$/caja 10 kilos"
$ws.Range("O282").Value = "China"
$ws.Range("P282").Value = 1800
$ws.Range("Q282").Value = 10
$ws.Range("R282").Value = "Hortaliza"

# Populate the second new row (becomes row 283) with fresh market data.
$ws.Range("A283").Value = 10
$ws.Range("B283").Value = "Vega Modelo de Temuco"
$ws.Range("C283").Value = "La Araucanía"
$ws.Range("D283").Value = 44463
$ws.Range("E283").Value = 9
$ws.Range("F283").Value = 100112003
$ws.Range("G283").Value = "Ajo"
$ws.Range("H283").Value = "Chino"
$ws.Range("I283").Value = "Primera"
$ws.Range("J283").Value = 80
$ws.Range("K283").Value = 19000
$ws.Range("L283").Value = 19000
$ws.Range("M283").Value = 19000
$ws.Range("N283").Value = "$/malla 10 kilos"
$ws.Range("O283").Value = "China"
$ws.Range("P283").Value = 1900
$ws.Range("Q283").Value = 10
$ws.Range("R283").Value = "Hortaliza"
